# Kayıt silindi: 11169351
# The record with Kayıt No 11169351 (row 796 on "Kayitlar", row 257 on
# "Merkez İlçe") is removed; every row below it shifts up by one on both
# sheets, shrinking each sheet's used range by one row.

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows(796).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows(257).Delete()
